# Rebuild the age_group summary table: header labels switch from the
# anaemia/no-anaemia wording to TRUE/FALSE wording, and the count/percent
# figures are recomputed (rows 6 & 7 additionally swap which age band they
# describe, with "35 to 39, years" becoming "40 to 44 years" in row 6 and
# "40 to 44 years" becoming "35 to 39 years" in row 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("B1").Value = "FALSE_count"
$ws.Range("C1").Value = "FALSE_percent"
$ws.Range("D1").Value = "TRUE_count"
$ws.Range("E1").Value = "TRUE_percent"

# --- Row 2: 15 to 19 years ---
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 16.66666666666666
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 24.32432432432433

# --- Row 3: 20 to 24 years ---
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 16.66666666666666
$ws.Range("D3").Value = 18
$ws.Range("E3").Value = 24.32432432432433

# --- Row 4: 25 to 29 years ---
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 38.88888888888889
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 24.32432432432433

# --- Row 5: 30 to 34 years ---
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 22.22222222222222
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 10.81081081081081

# --- Row 6: was "35 to 39, years", becomes "40 to 44 years" ---
$ws.Range("A6").Value = "40 to 44 years"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 5.555555555555555
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 5.405405405405405
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 7.692307692307693

# --- Row 7: was "40 to 44 years", becomes "35 to 39 years" ---
$ws.Range("A7").Value = "35 to 39 years"
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = 10.81081081081081
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.564102564102564
